$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-11 (Language, Count, Last User)
$data = @(
    @("de", 6,  "schaumburgernachrichten"),
    @("en", 19, "rff"),
    @("fi", 2,  "MinnaRuokonen"),
    @("pt", 1,  "g1_globo"),
    @("th", 1,  "ogataquotes_th"),
    @("ja", 4,  "osame120"),
    @("es", 1,  "MelchorRuizCope"),
    @("",   4,  "elysion"),
    @("tr", 1,  "tramboline"),
    @("zh", 1,  "tiefraum")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
